$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 467-468),
# pushing the existing rows 467-476 down to 469-478.
$ws.Rows("467:468").Insert()

# New row 467: Ají, Americana (o), Primera - week of 2023-11-09
$ws.Cells.Item(467,1).Value = 2
$ws.Cells.Item(467,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(467,3).Value = "Coquimbo"
$ws.Cells.Item(467,4).Value = 45239
$ws.Cells.Item(467,5).Value = 4
$ws.Cells.Item(467,6).Value = 100112021
$ws.Cells.Item(467,7).Value = "Ají"
$ws.Cells.Item(467,8).Value = "Americana (o)"
$ws.Cells.Item(467,9).Value = "Primera"
$ws.Cells.Item(467,10).Value = 160
$ws.Cells.Item(467,11).Value = 28000
$ws.Cells.Item(467,12).Value = 33000
$ws.Cells.Item(467,13).Value = 30500
$ws.Cells.Item(467,14).Value = "$/caja 25 kilos"
$ws.Cells.Item(467,15).Value = "Provincia de Limarí"
$ws.Cells.Item(467,16).Value = 1220
$ws.Cells.Item(467,17).Value = 25
$ws.Cells.Item(467,18).Value = "Hortaliza"

# New row 468: Ají, Inferno, Primera - week of 2023-11-09
$ws.Cells.Item(468,1).Value = 2
$ws.Cells.Item(468,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(468,3).Value = "Coquimbo"
$ws.Cells.Item(468,4).Value = 45239
$ws.Cells.Item(468,5).Value = 4
$ws.Cells.Item(468,6).Value = 100112021
$ws.Cells.Item(468,7).Value = "Ají"
$ws.Cells.Item(468,8).Value = "Inferno"
$ws.Cells.Item(468,9).Value = "Primera"
$ws.Cells.Item(468,10).Value = 160
$ws.Cells.Item(468,11).Value = 21000
$ws.Cells.Item(468,12).Value = 23000
$ws.Cells.Item(468,13).Value = 22000
$ws.Cells.Item(468,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(468,15).Value = "Provincia de Limarí"
$ws.Cells.Item(468,16).Value = 1467
$ws.Cells.Item(468,17).Value = 15
$ws.Cells.Item(468,18).Value = "Hortaliza"
